$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New "Route7" column (column H) added alongside the existing Route1..Route6 columns.
$ws.Range("H1").Value = "Route7"
$ws.Range("H2").Value = "PAD,BRI"
$ws.Range("H3").Value = "BRI,PAD"

# These rows hold values that look numeric ("0900,1000"); prefix with an
# apostrophe so Excel stores them as quoted text (quotePrefix style),
# matching the other data rows in the sheet.
$ws.Range("H4").Value = "'0900,1000"
$ws.Range("H5").Value = "'0900,1000"
$ws.Range("H6").Value = "'0900,1000"
$ws.Range("H7").Value = "'0900,1000"
$ws.Range("H8").Value = "'0900,1000"
$ws.Range("H9").Value = "'0900,1000"

# Move the active selection like in the target workbook.
$ws.Range("I4").Select()
